$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.517.75"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3
$ws.Range("D3").Value = "1.914.24"
$ws.Range("E3").Value = "  +2.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'315.23"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("D7").Value = "'0.5163"
$ws.Range("E7").Value = "  +3.57%  "

# Row 8
$ws.Range("D8").Value = "'0.3975"
$ws.Range("E8").Value = "  +0.19%  "

# Row 9
$ws.Range("D9").Value = "'0.09855"
$ws.Range("E9").Value = "  -2.44%  "

# Row 10
$ws.Range("D10").Value = "'1.151"
$ws.Range("E10").Value = "  +2.98%  "

# Row 11
$ws.Range("E11").Value = "  +2.54%  "

# Row 12
$ws.Range("D12").Value = "'6.530"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13
$ws.Range("D13").Value = "'21.16"
$ws.Range("E13").Value = "  +1.76%  "

# Row 14
$ws.Range("D14").Value = "1.921.34"
$ws.Range("E14").Value = "  +3.88%  "

# Row 15
$ws.Range("D15").Value = "'7.469"
$ws.Range("E15").Value = "  +1.13%  "

# Row 16
$ws.Range("E16").Value = "  +0.21%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001138"
$ws.Range("E17").Value = "  -1.11%  "

# Row 18
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'94.59"
$ws.Range("E18").Value = "  +1.18%  "

# Row 19
$ws.Range("D19").Value = "'0.06659"
$ws.Range("E19").Value = "  -0.02%  "

# Row 20
$ws.Range("D20").Value = "'18.23"
$ws.Range("E20").Value = "  +5.03%  "

# Row 21
$ws.Range("E21").Value = "  -0.02%  "

# Row 22
$ws.Range("D22").Value = "'6.314"
$ws.Range("E22").Value = "  +4.39%  "

# Row 23
$ws.Range("D23").Value = "28.572.90"
$ws.Range("E23").Value = "  +0.58%  "

# Row 24
$ws.Range("D24").Value = "'11.49"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
$ws.Range("D25").Value = "'2.323"
$ws.Range("E25").Value = "  +3.29%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.682"
$ws.Range("E26").Value = "  +8.60%  "

# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.130.17"
$ws.Range("E27").Value = "  +3.45%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.27"
$ws.Range("E28").Value = "  +0.79%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'157.63"
$ws.Range("E29").Value = "  +0.09%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'129.26"
$ws.Range("E30").Value = "  +1.25%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.112"
$ws.Range("E31").Value = "  +5.53%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1076"
$ws.Range("E32").Value = "  +2.30%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.744"
$ws.Range("E33").Value = "  +1.93%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.634"
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.888"
$ws.Range("E35").Value = "  +8.82%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06763"
$ws.Range("E36").Value = "  -0.46%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02438"
$ws.Range("E37").Value = "  +2.67%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.274"
$ws.Range("E38").Value = "  +5.42%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2218"
$ws.Range("E39").Value = "  +2.62%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.79"
$ws.Range("E40").Value = "  +2.60%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6480"
$ws.Range("E41").Value = "  +3.28%  "

# Row 42
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "'5.093"
$ws.Range("E42").Value = "  +1.31%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.188"
$ws.Range("E43").Value = "  +0.74%  "

# Row 44
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = "  +0.21%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.60"
$ws.Range("E45").Value = "  +1.86%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6100"
$ws.Range("E46").Value = "  +2.02%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.766"
$ws.Range("E47").Value = "  +1.98%  "

# Row 48
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.289"
$ws.Range("E48").Value = "  +0.82%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.050"
$ws.Range("E49").Value = "  +4.26%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'124.76"
$ws.Range("E50").Value = "  +0.11%  "

# Row 51
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.204"
$ws.Range("E51").Value = "  +1.17%  "
